# Sync attendance_reports: recorded-by swap + session 13/14 status swap

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#    Every G-column "Recorded By" cell that lists both the user and the
#    System account gets the order swapped.
# ---------------------------------------------------------------------
$gCells = @(8, 9, 10, 34, 35, 36, 60, 61, 62, 86, 87, 88, 112, 113, 114, `
            138, 139, 140, 164, 167, 191, 194, 218, 221, 245, 248, 272, `
            275, 299, 302)

foreach ($r in $gCells) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}

# ---------------------------------------------------------------------
# 2) Session 13 / Session 14 rows had their "Recorded"/"Not Recorded"
#    status (fill colour + Recorded By/Students/Status columns) on the
#    wrong row; swap each pair back.
# ---------------------------------------------------------------------
$rowPairs = @(
    @(13, 14),
    @(39, 40),
    @(65, 66),
    @(91, 92),
    @(117, 118),
    @(143, 144)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("A" + $r1 + ":I" + $r1)
    $range2 = $ws.Range("A" + $r2 + ":I" + $r2)

    # Swap the row fill colour (Recorded = green style, Not Recorded = pink style)
    $color1 = $range1.Interior.Color
    $color2 = $range2.Interior.Color
    $range1.Interior.Color = $color2
    $range2.Interior.Color = $color1

    # Swap the "Recorded By" / "Students" / "Status" values (G:I)
    $g1 = $ws.Cells.Item($r1, 7).Value2
    $h1 = $ws.Cells.Item($r1, 8).Value2
    $i1 = $ws.Cells.Item($r1, 9).Value2

    $g2 = $ws.Cells.Item($r2, 7).Value2
    $h2 = $ws.Cells.Item($r2, 8).Value2
    $i2 = $ws.Cells.Item($r2, 9).Value2

    $ws.Cells.Item($r1, 7).Value2 = $g2
    $ws.Cells.Item($r1, 8).Value2 = $h2
    $ws.Cells.Item($r1, 9).Value2 = $i2

    $ws.Cells.Item($r2, 7).Value2 = $g1
    $ws.Cells.Item($r2, 8).Value2 = $h1
    $ws.Cells.Item($r2, 9).Value2 = $i1
}
